# duplication & similarity checking completed
# Duplicate the "Which of these is a mammal?" question row on the
# MultipleChoice sheet, with a slightly reworded question text, to
# exercise the duplicate/similarity checker.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MultipleChoice")

# Append the duplicated question as a new row (row 9) — same answers,
# category and image as row 4, but reworded question text.
$ws.Range("A9").Value = "Which one is a mammal?"
$ws.Range("B9").Value = "Shark"
$ws.Range("C9").Value = "Frog"
$ws.Range("D9").Value = "Dolphin"
$ws.Range("E9").Value = "Eagle"
$ws.Range("F9").Value = "Lizard"
$ws.Range("G9").Value = "c"
$ws.Range("H9").Value = "dolphin.png"
$ws.Range("I9").Value = "animal"

# Move focus to the MultipleChoice sheet and leave the selection where
# the author last clicked while reviewing the new row.
$ws.Activate()
$ws.Range("A11").Select()
